# Updated symbol list on Sat Jan  7 14:56:35 UTC 2023 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) values for the
# cryptocurrency rows on Sheet1, keeping each cell as text (matching the
# original inline-string storage) by using a leading apostrophe and then
# resetting the cell style back to "Normal" so no stray number-format /
# quote-prefix style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "'260.34"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = "'1.46%"
$ws.Cells.Item(2,5).Style = "Normal"
$ws.Cells.Item(3,4).Value = "'27.26"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = "'2.51%"
$ws.Cells.Item(3,5).Style = "Normal"
$ws.Cells.Item(4,4).Value = "'4.676"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = "'0.70%"
$ws.Cells.Item(4,5).Style = "Normal"
$ws.Cells.Item(5,4).Value = "'0.06171"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "'4.21%"
$ws.Cells.Item(5,5).Style = "Normal"
$ws.Cells.Item(6,4).Value = "'6.669"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = "'1.04%"
$ws.Cells.Item(6,5).Style = "Normal"
$ws.Cells.Item(7,4).Value = "'0.8510"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = "'-0.52%"
$ws.Cells.Item(7,5).Style = "Normal"
$ws.Cells.Item(8,4).Value = "'0.9138"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "'0.06%"
$ws.Cells.Item(8,5).Style = "Normal"
$ws.Cells.Item(9,4).Value = "'0.1408"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = "'2.26%"
$ws.Cells.Item(9,5).Style = "Normal"
$ws.Cells.Item(10,4).Value = "'0.04849"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "'8.20%"
$ws.Cells.Item(10,5).Style = "Normal"
$ws.Cells.Item(11,4).Value = "'0.07092"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "'1.36%"
$ws.Cells.Item(11,5).Style = "Normal"
$ws.Cells.Item(12,5).Value = "'2.36%"
$ws.Cells.Item(12,5).Style = "Normal"
$ws.Cells.Item(13,4).Value = "'0.09058"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = "'-0.49%"
$ws.Cells.Item(13,5).Style = "Normal"
$ws.Cells.Item(14,4).Value = "'0.001542"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "'1.00%"
$ws.Cells.Item(14,5).Style = "Normal"
$ws.Cells.Item(15,4).Value = "'0.0006187"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = "'2.60%"
$ws.Cells.Item(15,5).Style = "Normal"
$ws.Cells.Item(16,4).Value = "'0.006106"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = "'0.77%"
$ws.Cells.Item(16,5).Style = "Normal"
$ws.Cells.Item(17,5).Value = "'-0.34%"
$ws.Cells.Item(17,5).Style = "Normal"
$ws.Cells.Item(18,5).Value = "'0.63%"
$ws.Cells.Item(18,5).Style = "Normal"
$ws.Cells.Item(19,4).Value = "'2.180"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "'0.57%"
$ws.Cells.Item(19,5).Style = "Normal"
$ws.Cells.Item(20,5).Value = "'-0.01%"
$ws.Cells.Item(20,5).Style = "Normal"
$ws.Cells.Item(21,4).Value = "'0.1299"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "'0.85%"
$ws.Cells.Item(21,5).Style = "Normal"
$ws.Cells.Item(22,4).Value = "'4.096"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = "'5.74%"
$ws.Cells.Item(22,5).Style = "Normal"
$ws.Cells.Item(23,4).Value = "'0.04245"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = "'0.86%"
$ws.Cells.Item(23,5).Style = "Normal"
$ws.Cells.Item(24,5).Value = "'0.09%"
$ws.Cells.Item(24,5).Style = "Normal"
$ws.Cells.Item(25,4).Value = "'0.003796"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "'-17.66%"
$ws.Cells.Item(25,5).Style = "Normal"
$ws.Cells.Item(26,5).Value = "'0.05%"
$ws.Cells.Item(26,5).Style = "Normal"
$ws.Cells.Item(27,4).Value = "'0.0001575"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "'-8.19%"
$ws.Cells.Item(27,5).Style = "Normal"
$ws.Cells.Item(40,4).Value = "'0.03874"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "'2.04%"
$ws.Cells.Item(40,5).Style = "Normal"
$ws.Cells.Item(41,4).Value = "'0.1114"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "'1.46%"
$ws.Cells.Item(41,5).Style = "Normal"
$ws.Cells.Item(42,4).Value = "'0.004081"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = "'-34.67%"
$ws.Cells.Item(42,5).Style = "Normal"
$ws.Cells.Item(43,5).Value = "'13.84%"
$ws.Cells.Item(43,5).Style = "Normal"
$ws.Cells.Item(44,4).Value = "'0.002196"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "'-4.80%"
$ws.Cells.Item(44,5).Style = "Normal"
$ws.Cells.Item(45,4).Value = "'0.00005162"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "'1.06%"
$ws.Cells.Item(45,5).Style = "Normal"
$ws.Cells.Item(46,5).Value = "'-0.01%"
$ws.Cells.Item(46,5).Style = "Normal"
$ws.Cells.Item(47,5).Value = "'7.99%"
$ws.Cells.Item(47,5).Style = "Normal"
$ws.Cells.Item(49,4).Value = "'0.00002100"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = "'-0.01%"
$ws.Cells.Item(49,5).Style = "Normal"
$ws.Cells.Item(50,4).Value = "'0.0002000"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "'-0.01%"
$ws.Cells.Item(50,5).Style = "Normal"
